# Applies the renaming of sheets and the selection/active-tab changes
# described by the commit: sheet names simplified via a "RegEx rid off"
# pass, and the active sheet / selections updated to reflect where the
# author was last working.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (input_* / *_names -> shorter names) -------------------
$wb.Worksheets.Item(1).Name = "stoich_coeff"
$wb.Worksheets.Item(2).Name = "constants_log10"
$wb.Worksheets.Item(3).Name = "concentra"
$wb.Worksheets.Item(4).Name = "component"

# --- Update per-sheet selections -------------------------------------------
# Sheet 1 (stoich_coeff): selection moves from J11 to G31
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("G31").Select()

# Sheet 3 (concentra): selection collapses from F2:F13 to H8 and it is no
# longer the tab that is selected
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("H8").Select()

# Sheet 4 (component): becomes the active/selected tab, with selection I18
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
$ws4.Range("I18").Select()
